$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the RUNS and WICKETS columns still held numbers left over from a
# previous tournament. Reset every player's stats back to zero so the
# standings sheet starts clean for the new tournament.
for ($r = 2; $r -le 89; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
}

# Widen the player-name column so the longer names are fully visible.
$ws.Columns.Item(1).ColumnWidth = 32.14

# Restore the cursor/selection position that was active when the fix was made.
$ws.Range("F73").Select() | Out-Null
